$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.993.42"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.578.38"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "2.595.54"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.79"
$ws.Range("D10").ClearFormats()
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "3.034.49"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "58.004.15"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "2.599.81"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -5.06%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "2.693.01"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "0.0₃0744"
$ws.Range("E30").Value = "  -6.60%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.75"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.66"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.867"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "271.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0953"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0523"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.971.04"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.50%  "
